# Trade #135 closed at 2026-02-18 00:41:21 - unknown UNKNOWN +0.000%
#
# Applies:
#  - Summary sheet roll-up metrics refresh
#  - Strategy Status roll-up refresh for MarketMaking
#  - "All Trades" row 164 (MarketMaking trade, entered 00:37:35) closes out
#  - Two brand-new OPEN trades get logged (momentum #192, MarketMaking #193)
#    into "All Trades" plus their respective per-strategy sheets
#  - The per-strategy "MarketMaking" sheet row 64 mirrors the same close

$wb = $excel.ActiveWorkbook

# Helper: write a literal text value into a cell without Excel's
# autoconvert-to-date/number kicking in (e.g. "2026-02-18" -> date serial).
# Force the cell to Text format, assign, then drop the format back to
# Normal/General so no stray style survives in the saved file.
function Set-LiteralText {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.11   # Current Capital
$summary.Range("B4").Value = 0.22      # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 163       # Total Trades
$summary.Range("B7").Value = 74        # Winning Trades
$summary.Range("B9").Value = 45.4      # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Cells.Item(6, 3).Value = 99.31                 # Capital
$status.Cells.Item(6, 4).Value = 63                    # Trades
$status.Cells.Item(6, 5).Value = -0.5                  # P&L $
$status.Cells.Item(6, 6).Value = -0.6899999999999999   # P&L %
$status.Cells.Item(6, 7).Value = 47.62                 # Win Rate %

# ---------------------------------------------------------------------
# "All Trades" sheet
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Row 164: close out the existing MarketMaking trade (#163)
$allTrades.Cells.Item(164, 7).Value = 0.83              # Exit Price
$allTrades.Cells.Item(164, 8).Value = "CLOSED"          # Status
$allTrades.Cells.Item(164, 9).Value = 2.4691            # P&L %
$allTrades.Cells.Item(164, 10).Value = 0.02             # P&L $
$allTrades.Cells.Item(164, 11).Value = 99.31            # Capital After
$allTrades.Cells.Item(164, 12).Value = "early_exit"     # Exit Reason
$allTrades.Cells.Item(164, 13).Value = 0.18             # Duration (min)

# Row 193: new OPEN trade (#192, momentum)
$allTrades.Cells.Item(193, 1).Value = 192
Set-LiteralText $allTrades.Cells.Item(193, 2) "2026-02-18"
$allTrades.Cells.Item(193, 3).Value = "00:41:15"
$allTrades.Cells.Item(193, 4).Value = "momentum"
$allTrades.Cells.Item(193, 5).Value = "DOWN"
$allTrades.Cells.Item(193, 6).Value = 0.8100000000000001
$allTrades.Cells.Item(193, 8).Value = "OPEN"
$allTrades.Cells.Item(193, 9).Value = 0
$allTrades.Cells.Item(193, 10).Value = 0
$allTrades.Cells.Item(193, 11).Value = 99.15712996249174
$allTrades.Cells.Item(193, 13).Value = 0
$allTrades.Cells.Item(193, 14).Value = 0
$allTrades.Cells.Item(193, 15).Value = 0
$allTrades.Cells.Item(193, 16).Value = 0.9
$allTrades.Cells.Item(193, 17).Value = "Downward momentum: -45.109% over 10 samples"

# Row 194: new OPEN trade (#193, MarketMaking)
$allTrades.Cells.Item(194, 1).Value = 193
Set-LiteralText $allTrades.Cells.Item(194, 2) "2026-02-18"
$allTrades.Cells.Item(194, 3).Value = "00:41:16"
$allTrades.Cells.Item(194, 4).Value = "MarketMaking"
$allTrades.Cells.Item(194, 5).Value = "DOWN"
$allTrades.Cells.Item(194, 6).Value = 0.82
$allTrades.Cells.Item(194, 8).Value = "OPEN"
$allTrades.Cells.Item(194, 9).Value = 0
$allTrades.Cells.Item(194, 10).Value = 0
$allTrades.Cells.Item(194, 11).Value = 99.28858346467945
$allTrades.Cells.Item(194, 13).Value = 0
$allTrades.Cells.Item(194, 14).Value = 0
$allTrades.Cells.Item(194, 15).Value = 0
$allTrades.Cells.Item(194, 16).Value = 0.6
$allTrades.Cells.Item(194, 17).Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# "momentum" sheet - new OPEN trade (#192)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(52, 1).Value = 192
Set-LiteralText $momentum.Cells.Item(52, 2) "2026-02-18"
$momentum.Cells.Item(52, 3).Value = "00:41:15"
$momentum.Cells.Item(52, 4).Value = "momentum"
$momentum.Cells.Item(52, 5).Value = "DOWN"
$momentum.Cells.Item(52, 6).Value = 0.8100000000000001
$momentum.Cells.Item(52, 8).Value = "OPEN"
$momentum.Cells.Item(52, 9).Value = 0
$momentum.Cells.Item(52, 10).Value = 0
$momentum.Cells.Item(52, 11).Value = 99.15712996249174
$momentum.Cells.Item(52, 12).Value = 0
$momentum.Cells.Item(52, 13).Value = 0
$momentum.Cells.Item(52, 14).Value = 0.9
$momentum.Cells.Item(52, 15).Value = "Downward momentum: -45.109% over 10 samples"
$momentum.Cells.Item(52, 17).Value = 0

# ---------------------------------------------------------------------
# "MarketMaking" sheet
# ---------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# Row 64: close out the existing trade (#163) - mirrors All Trades row 164
$marketMaking.Cells.Item(64, 7).Value = 0.83            # Exit Price
$marketMaking.Cells.Item(64, 8).Value = "CLOSED"        # Status
$marketMaking.Cells.Item(64, 9).Value = 2.4691          # P&L %
$marketMaking.Cells.Item(64, 10).Value = 0.02           # P&L $
$marketMaking.Cells.Item(64, 11).Value = 99.31          # Capital After
$marketMaking.Cells.Item(64, 16).Value = "early_exit"   # Exit Reason
$marketMaking.Cells.Item(64, 17).Value = 0.18           # Duration (min)

# Row 82: new OPEN trade (#193)
$marketMaking.Cells.Item(82, 1).Value = 193
Set-LiteralText $marketMaking.Cells.Item(82, 2) "2026-02-18"
$marketMaking.Cells.Item(82, 3).Value = "00:41:16"
$marketMaking.Cells.Item(82, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(82, 5).Value = "DOWN"
$marketMaking.Cells.Item(82, 6).Value = 0.82
$marketMaking.Cells.Item(82, 8).Value = "OPEN"
$marketMaking.Cells.Item(82, 9).Value = 0
$marketMaking.Cells.Item(82, 10).Value = 0
$marketMaking.Cells.Item(82, 11).Value = 99.28858346467945
$marketMaking.Cells.Item(82, 12).Value = 0
$marketMaking.Cells.Item(82, 13).Value = 0
$marketMaking.Cells.Item(82, 14).Value = 0.6
$marketMaking.Cells.Item(82, 15).Value = "Normal spread capture: 198 bps"
$marketMaking.Cells.Item(82, 17).Value = 0
